# RDCC-3540 Upload Template Improvements
# - Rename "Case Worker Data" sheet to "Staff Data"
# - Rename header columns on that sheet from "Area of WorkN" / "Area of WorkN ID"
#   to "ServiceN" / "ServiceN ID" (N = 1..8)
# - Update the sheet's selection/scroll position

$wb = $excel.ActiveWorkbook

# --- Rename the "Case Worker Data" sheet to "Staff Data" ---
$ws = $wb.Worksheets.Item("Case Worker Data")
$ws.Name = "Staff Data"

# --- Rename the "Area of WorkN" / "Area of WorkN ID" header pairs to "ServiceN" / "ServiceN ID" ---
# Row 1 columns: M/O/Q/S/U/W/Y/AA hold "Area of WorkN" (N=1..8)
#                N/P/R/T/V/X/Z/AB hold "Area of WorkN ID" (N=1..8)
# Written in N,O,P,Q,R,S,T,U,V,W,X,Y,Z,AA,AB,M order (M last) to match the
# original authoring order of the new shared strings.
$ws.Range("N1").Value = "Service1 ID"
$ws.Range("O1").Value = "Service2"
$ws.Range("P1").Value = "Service2 ID"
$ws.Range("Q1").Value = "Service3"
$ws.Range("R1").Value = "Service3 ID"
$ws.Range("S1").Value = "Service4"
$ws.Range("T1").Value = "Service4 ID"
$ws.Range("U1").Value = "Service5"
$ws.Range("V1").Value = "Service5 ID"
$ws.Range("W1").Value = "Service6"
$ws.Range("X1").Value = "Service6 ID"
$ws.Range("Y1").Value = "Service7"
$ws.Range("Z1").Value = "Service7 ID"
$ws.Range("AA1").Value = "Service8"
$ws.Range("AB1").Value = "Service8 ID"
$ws.Range("M1").Value = "Service1"

# --- Update view: scroll position and selection ---
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 8
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("M13").Select()

Write-Output "Renamed sheet and updated Service headers"
